$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3429493964259791
$ws.Cells.Item(2, 3).Value = 0.06965172002792031
$ws.Cells.Item(2, 4).Value = 0.0261353511007556
$ws.Cells.Item(2, 5).Value = 0.4172994359769149
$ws.Cells.Item(2, 6).Value = 0.6289244698072523
$ws.Cells.Item(2, 9).Value = 0.472506710197468
$ws.Cells.Item(2, 11).Value = 0.3790598273342312
$ws.Cells.Item(2, 15).Value = 2.123885881344364

$ws.Cells.Item(3, 2).Value = 0.3009604902049148
$ws.Cells.Item(3, 3).Value = 0.06153662247182012
$ws.Cells.Item(3, 4).Value = 0.02409310344395266
$ws.Cells.Item(3, 5).Value = 0.3641208024745879
$ws.Cells.Item(3, 6).Value = 0.6278367374468559
$ws.Cells.Item(3, 9).Value = 0.4777763309126399
$ws.Cells.Item(3, 11).Value = 0.3309631340928831
$ws.Cells.Item(3, 15).Value = 2.133867296146576

$ws.Cells.Item(4, 2).Value = 0.2751380727284243
$ws.Cells.Item(4, 3).Value = 0.05652751646485399
$ws.Cells.Item(4, 4).Value = 0.02282952962735152
$ws.Cells.Item(4, 5).Value = 0.3315522843216456
$ws.Cells.Item(4, 6).Value = 0.6276358551096237
$ws.Cells.Item(4, 9).Value = 0.4813320661462512
$ws.Cells.Item(4, 11).Value = 0.3013555354349648
$ws.Cells.Item(4, 15).Value = 2.141531551795708

$ws.Cells.Item(5, 2).Value = 0.2646054256214541
$ws.Cells.Item(5, 3).Value = 0.05447967130379539
$ws.Cells.Item(5, 4).Value = 0.02231222978161185
$ws.Cells.Item(5, 5).Value = 0.3182997360578526
$ws.Cells.Item(5, 6).Value = 0.6276712106583346
$ws.Cells.Item(5, 9).Value = 0.4828614232459891
$ws.Cells.Item(5, 11).Value = 0.2892715992569777
$ws.Cells.Item(5, 15).Value = 2.145040146667014

$ws.Cells.Item(6, 2).Value = 0.2628559110286233
$ws.Cells.Item(6, 3).Value = 0.05413923129175657
$ws.Cells.Item(6, 4).Value = 0.02222618969229018
$ws.Cells.Item(6, 5).Value = 0.316100285483131
$ws.Cells.Item(6, 6).Value = 0.6276841548369489
$ws.Cells.Item(6, 9).Value = 0.4831202215763213
$ws.Cells.Item(6, 11).Value = 0.287263962931803
$ws.Cells.Item(6, 15).Value = 2.145645996966479

$ws.Cells.Item(7, 2).Value = 0.2749960647591365
$ws.Cells.Item(7, 3).Value = 0.05649992514496205
$ws.Cells.Item(7, 4).Value = 0.02282256274047967
$ws.Cells.Item(7, 5).Value = 0.3313734793941308
$ws.Cells.Item(7, 6).Value = 0.6276358576039129
$ws.Cells.Item(7, 9).Value = 0.4813523663823851
$ws.Cells.Item(7, 11).Value = 0.3011926418602684
$ws.Cells.Item(7, 15).Value = 2.141577310738967

$ws.Cells.Item(8, 2).Value = 0.3284804428746781
$ws.Cells.Item(8, 3).Value = 0.0668591391214477
$ws.Cells.Item(8, 4).Value = 0.02543320349975886
$ws.Cells.Item(8, 5).Value = 0.3989451711657068
$ws.Cells.Item(8, 6).Value = 0.628452366240829
$ws.Cells.Item(8, 9).Value = 0.4742571391747639
$ws.Cells.Item(8, 11).Value = 0.3624921648033137
$ws.Cells.Item(8, 15).Value = 2.127008269994889

$ws.Cells.Item(9, 2).Value = 0.4330196250673737
$ws.Cells.Item(9, 3).Value = 0.08696300939867285
$ws.Cells.Item(9, 4).Value = 0.03047489210481302
$ws.Cells.Item(9, 5).Value = 0.5321920092804504
$ws.Cells.Item(9, 6).Value = 0.6337695595010473
$ws.Cells.Item(9, 9).Value = 0.4628902603821636
$ws.Cells.Item(9, 11).Value = 0.4820801107834711
$ws.Cells.Item(9, 15).Value = 2.11066090298516

$ws.Cells.Item(10, 2).Value = 0.5095990029419966
$ws.Cells.Item(10, 3).Value = 0.1016050693608008
$ws.Cells.Item(10, 4).Value = 0.03413008150782559
$ws.Cells.Item(10, 5).Value = 0.6306550029007383
$ws.Cells.Item(10, 6).Value = 0.639957876948877
$ws.Cells.Item(10, 9).Value = 0.4561006330277664
$ws.Cells.Item(10, 11).Value = 0.5695492469237422
$ws.Cells.Item(10, 15).Value = 2.106157102289757

$ws.Cells.Item(11, 2).Value = 0.5443850744516681
$ws.Cells.Item(11, 3).Value = 0.1082383929676212
$ws.Cells.Item(11, 4).Value = 0.03578199605681931
$ws.Cells.Item(11, 5).Value = 0.6755970522737158
$ws.Cells.Item(11, 6).Value = 0.6432722012526071
$ws.Cells.Item(11, 9).Value = 0.4533528994753198
$ws.Cells.Item(11, 11).Value = 0.6092537054501292
$ws.Cells.Item(11, 15).Value = 2.105750366634055

$ws.Cells.Item(12, 2).Value = 0.5575500160009597
$ws.Cells.Item(12, 3).Value = 0.1107462947907436
$ws.Cells.Item(12, 4).Value = 0.03640594065365832
$ws.Cells.Item(12, 5).Value = 0.6926390146231398
$ws.Cells.Item(12, 6).Value = 0.6445993015316915
$ws.Cells.Item(12, 9).Value = 0.4523615882685732
$ws.Cells.Item(12, 11).Value = 0.624276038997948
$ws.Cells.Item(12, 15).Value = 2.105833419773148

$ws.Cells.Item(13, 2).Value = 0.5547150655617372
$ws.Cells.Item(13, 3).Value = 0.1102063516823364
$ws.Cells.Item(13, 4).Value = 0.03627163478750361
$ws.Cells.Item(13, 5).Value = 0.6889676557752011
$ws.Cells.Item(13, 6).Value = 0.6443102788797432
$ws.Cells.Item(13, 9).Value = 0.452572894090256
$ws.Cells.Item(13, 11).Value = 0.6210412925493358
$ws.Cells.Item(13, 15).Value = 2.105804973691619

$ws.Cells.Item(14, 2).Value = 0.5454683214058207
$ws.Cells.Item(14, 3).Value = 0.1084447998888436
$ws.Cells.Item(14, 4).Value = 0.03583336063678644
$ws.Cells.Item(14, 5).Value = 0.6769986262960401
$ws.Cells.Item(14, 6).Value = 0.6433799374226723
$ws.Cells.Item(14, 9).Value = 0.453270356746259
$ws.Cells.Item(14, 11).Value = 0.610489863266821
$ws.Cells.Item(14, 15).Value = 2.105752442689777

$ws.Cells.Item(15, 2).Value = 0.5398033927484107
$ws.Cells.Item(15, 3).Value = 0.1073652769509863
$ws.Cells.Item(15, 4).Value = 0.0355646955850375
$ws.Cells.Item(15, 5).Value = 0.6696703463099993
$ws.Cells.Item(15, 6).Value = 0.6428194652370678
$ws.Cells.Item(15, 9).Value = 0.4537039847020523
$ws.Cells.Item(15, 11).Value = 0.6040251138251449
$ws.Cells.Item(15, 15).Value = 2.105751168284826

$ws.Cells.Item(16, 2).Value = 0.5073245918154612
$ws.Cells.Item(16, 3).Value = 0.1011710099078869
$ws.Cells.Item(16, 4).Value = 0.03402190330888999
$ws.Cells.Item(16, 5).Value = 0.6277211131709066
$ws.Cells.Item(16, 6).Value = 0.6397513451047416
$ws.Cells.Item(16, 9).Value = 0.4562870786100106
$ws.Cells.Item(16, 11).Value = 0.5669526948046837
$ws.Cells.Item(16, 15).Value = 2.106216806973492

$ws.Cells.Item(17, 2).Value = 0.4873866239896927
$ws.Cells.Item(17, 3).Value = 0.097363967834184
$ws.Cells.Item(17, 4).Value = 0.03307264324255499
$ws.Cells.Item(17, 5).Value = 0.6020264356574216
$ws.Cells.Item(17, 6).Value = 0.6379971943105502
$ws.Cells.Item(17, 9).Value = 0.4579591603171949
$ws.Cells.Item(17, 11).Value = 0.5441876184391106
$ws.Cells.Item(17, 15).Value = 2.106923707552056

$ws.Cells.Item(18, 2).Value = 0.4759141336213872
$ws.Cells.Item(18, 3).Value = 0.09517167971392837
$ws.Cells.Item(18, 4).Value = 0.03252563509580852
$ws.Cells.Item(18, 5).Value = 0.5872616221732301
$ws.Cells.Item(18, 6).Value = 0.6370352279283651
$ws.Cells.Item(18, 9).Value = 0.4589529743305718
$ws.Cells.Item(18, 11).Value = 0.5310857162882598
$ws.Cells.Item(18, 15).Value = 2.107484800232072

$ws.Cells.Item(19, 2).Value = 0.4720289556434238
$ws.Cells.Item(19, 3).Value = 0.09442896736862849
$ws.Cells.Item(19, 4).Value = 0.0323402537526718
$ws.Cells.Item(19, 5).Value = 0.5822648774036026
$ws.Cells.Item(19, 6).Value = 0.636717583009542
$ws.Cells.Item(19, 9).Value = 0.4592949665982218
$ws.Cells.Item(19, 11).Value = 0.5266482749543684
$ws.Cells.Item(19, 15).Value = 2.107701283632878

$ws.Cells.Item(20, 2).Value = 0.4895095450594908
$ws.Cells.Item(20, 3).Value = 0.0977695011667663
$ws.Cells.Item(20, 4).Value = 0.03317379933214681
$ws.Cells.Item(20, 5).Value = 0.6047602078498784
$ws.Cells.Item(20, 6).Value = 0.6381790631378053
$ws.Cells.Item(20, 9).Value = 0.4577778432667827
$ws.Cells.Item(20, 11).Value = 0.5466118350788349
$ws.Cells.Item(20, 15).Value = 2.106832459512731

$ws.Cells.Item(21, 2).Value = 0.5481845301332271
$ws.Cells.Item(21, 3).Value = 0.1089623191431031
$ws.Cells.Item(21, 4).Value = 0.03596213610399701
$ws.Cells.Item(21, 5).Value = 0.6805135720823756
$ws.Cells.Item(21, 6).Value = 0.6436512441937623
$ws.Cells.Item(21, 9).Value = 0.4530641585026842
$ws.Cells.Item(21, 11).Value = 0.6135894279109664
$ws.Cells.Item(21, 15).Value = 2.105761430600836

$ws.Cells.Item(22, 2).Value = 0.586486369136594
$ws.Cells.Item(22, 3).Value = 0.1162541706304694
$ws.Cells.Item(22, 4).Value = 0.03777513364497764
$ws.Cells.Item(22, 5).Value = 0.7301598754764598
$ws.Cells.Item(22, 6).Value = 0.6476476368203947
$ws.Cells.Item(22, 9).Value = 0.4502703254260823
$ws.Cells.Item(22, 11).Value = 0.6572878066461669
$ws.Cells.Item(22, 15).Value = 2.106443748055881

$ws.Cells.Item(23, 2).Value = 0.5660483181394795
$ws.Cells.Item(23, 3).Value = 0.1123645208458299
$ws.Cells.Item(23, 4).Value = 0.03680837046060503
$ws.Cells.Item(23, 5).Value = 0.7036495962146603
$ws.Cells.Item(23, 6).Value = 0.6454761739146306
$ws.Cells.Item(23, 9).Value = 0.4517351437527459
$ws.Cells.Item(23, 11).Value = 0.633972238149056
$ws.Cells.Item(23, 15).Value = 2.105952788319257

$ws.Cells.Item(24, 2).Value = 0.4885498032003284
$ws.Cells.Item(24, 3).Value = 0.09758617066174224
$ws.Cells.Item(24, 4).Value = 0.03312807060387257
$ws.Cells.Item(24, 5).Value = 0.6035242465564465
$ws.Cells.Item(24, 6).Value = 0.6380966953545339
$ws.Cells.Item(24, 9).Value = 0.4578597154419271
$ws.Cells.Item(24, 11).Value = 0.5455158901696109
$ws.Cells.Item(24, 15).Value = 2.106873230967182

$ws.Cells.Item(25, 2).Value = 0.4047774681255589
$ws.Cells.Item(25, 3).Value = 0.08154691188582319
$ws.Cells.Item(25, 4).Value = 0.02911946681465594
$ws.Cells.Item(25, 5).Value = 0.4960541509791057
$ws.Cells.Item(25, 6).Value = 0.631931598301982
$ws.Cells.Item(25, 9).Value = 0.465691717167104
$ws.Cells.Item(25, 11).Value = 0.4497962708934722
$ws.Cells.Item(25, 15).Value = 2.113769074292236

